$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 803.61536
$ws.Range("I12").Value = 680.875
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 680.875
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -510.875
$ws.Range("N12").Value = -1340
$ws.Range("H58").Value = 425.66666
$ws.Range("I58").Value = 425.66666
$ws.Range("K58").Value = 1276.99998
$ws.Range("M58").Value = -1126.99998
$ws.Range("H106").Value = 3645.0908
$ws.Range("I106").Value = 3017.25
$ws.Range("K106").Value = 3017.25
$ws.Range("M106").Value = -2386.25
$ws.Range("H138").Value = 3631.2
$ws.Range("I138").Value = 4963
$ws.Range("J138").Value = 3029.742
$ws.Range("K138").Value = 14889
$ws.Range("L138").Value = 9089.226000000001
$ws.Range("M138").Value = -9749
$ws.Range("N138").Value = -19369.226
$ws.Range("H141").Value = 1558031.2
$ws.Range("I141").Value = 2334676.8
$ws.Range("K141").Value = 7004030.399999999
$ws.Range("M141").Value = -6998850.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7707.23
$ws.Range("I32").Value = 6377.446
$ws.Range("J32").Value = 22999.75
$ws.Range("K32").Value = 6377.446
$ws.Range("L32").Value = 22999.75
$ws.Range("M32").Value = -6090.446
$ws.Range("N32").Value = -23573.75
$ws.Range("H37").Value = 11600
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20546
$ws.Range("H63").Value = 2327.75
$ws.Range("I63").Value = 1602.5
$ws.Range("K63").Value = 1602.5
$ws.Range("M63").Value = -916.5
$ws.Range("H66").Value = 2327.75
$ws.Range("I66").Value = 1602.5
$ws.Range("K66").Value = 8012.5
$ws.Range("M66").Value = -4580.5
$ws.Range("H74").Value = 950.70734
$ws.Range("I74").Value = 775.5263
$ws.Range("J74").Value = 3169.6667
$ws.Range("K74").Value = 775.5263
$ws.Range("L74").Value = 3169.6667
$ws.Range("M74").Value = 98.47370000000001
$ws.Range("N74").Value = -4917.6667
$ws.Range("H77").Value = 950.70734
$ws.Range("I77").Value = 775.5263
$ws.Range("J77").Value = 3169.6667
$ws.Range("K77").Value = 3877.6315
$ws.Range("L77").Value = 15848.3335
$ws.Range("M77").Value = 490.3685
$ws.Range("N77").Value = -24584.3335
$ws.Range("H111").Value = 70000
$ws.Range("J111").Value = 70000
$ws.Range("L111").Value = 70000
$ws.Range("N111").Value = -78180
$ws.Range("H132").Value = 1855.0145
$ws.Range("I132").Value = 1464.7441
$ws.Range("J132").Value = 2500.4614
$ws.Range("K132").Value = 4394.2323
$ws.Range("L132").Value = 7501.3842
$ws.Range("M132").Value = -1864.2323
$ws.Range("N132").Value = -12561.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2439.8
$ws.Range("I99").Value = 2439.8
$ws.Range("K99").Value = 2439.8
$ws.Range("M99").Value = -941.8000000000002
$ws.Range("H107").Value = 2007.0769
$ws.Range("I107").Value = 2116
$ws.Range("K107").Value = 2116
$ws.Range("M107").Value = -196

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1358.5714
$ws.Range("J22").Value = 1847.5
$ws.Range("L22").Value = 1847.5
$ws.Range("N22").Value = -2547.5
$ws.Range("H31").Value = 2767.125
$ws.Range("I31").Value = 2315.6843
$ws.Range("J31").Value = 4482.6
$ws.Range("K31").Value = 2315.6843
$ws.Range("L31").Value = 4482.6
$ws.Range("M31").Value = -2020.6843
$ws.Range("N31").Value = -5072.6
$ws.Range("H34").Value = 2767.125
$ws.Range("I34").Value = 2315.6843
$ws.Range("J34").Value = 4482.6
$ws.Range("K34").Value = 2315.6843
$ws.Range("L34").Value = 4482.6
$ws.Range("M34").Value = -2113.6843
$ws.Range("N34").Value = -4886.6
$ws.Range("H58").Value = 854131.7
$ws.Range("I58").Value = 3346104.5
$ws.Range("K58").Value = 3346104.5
$ws.Range("M58").Value = -3345901.5
$ws.Range("H134").Value = 1600.9048
$ws.Range("I134").Value = 1469.9375
$ws.Range("J134").Value = 2020
$ws.Range("K134").Value = 4409.8125
$ws.Range("L134").Value = 6060
$ws.Range("M134").Value = -1874.8125
$ws.Range("N134").Value = -11130
$ws.Range("H136").Value = 854131.7
$ws.Range("I136").Value = 3346104.5
$ws.Range("K136").Value = 10038313.5
$ws.Range("M136").Value = -10035763.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 498
$ws.Range("I92").Value = 498
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1494
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H137").Value = 6678.4443
$ws.Range("J137").Value = 7031.467
$ws.Range("L137").Value = 21094.401
$ws.Range("N137").Value = -31294.401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 942299.25
$ws.Range("I132").Value = 1245271.6
$ws.Range("J132").Value = 3084.8
$ws.Range("K132").Value = 3735814.8
$ws.Range("L132").Value = 9254.400000000001
$ws.Range("M132").Value = -3733284.8
$ws.Range("N132").Value = -14314.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 6500
$ws.Range("I13").Value = 4000
$ws.Range("J13").Value = 9000
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 9000
$ws.Range("M13").Value = -3860
$ws.Range("N13").Value = -9280
$ws.Range("H40").Value = 17031.846
$ws.Range("I40").Value = 20189
$ws.Range("J40").Value = 11980.4
$ws.Range("K40").Value = 20189
$ws.Range("L40").Value = 11980.4
$ws.Range("M40").Value = -20053
$ws.Range("N40").Value = -12252.4
$ws.Range("H46").Value = 2033.8
$ws.Range("J46").Value = 2269.6
$ws.Range("L46").Value = 2269.6
$ws.Range("N46").Value = -2645.6
$ws.Range("H55").Value = 8334112.5
$ws.Range("J55").Value = 837.5
$ws.Range("L55").Value = 837.5
$ws.Range("N55").Value = -1183.5
$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -22350
$ws.Range("H136").Value = 1309.1
$ws.Range("I136").Value = 948.25
$ws.Range("J136").Value = 2752.5
$ws.Range("K136").Value = 2844.75
$ws.Range("L136").Value = 8257.5
$ws.Range("M136").Value = -294.75
$ws.Range("N136").Value = -13357.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 13249
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H126").Value = 2916.1428
$ws.Range("I126").Value = 2486
$ws.Range("K126").Value = 7458
$ws.Range("M126").Value = -4988
$ws.Range("H132").Value = 1762.125
$ws.Range("I132").Value = 1604.9
$ws.Range("J132").Value = 2548.25
$ws.Range("K132").Value = 4814.700000000001
$ws.Range("L132").Value = 7644.75
$ws.Range("M132").Value = -2284.700000000001
$ws.Range("N132").Value = -12704.75
